$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "61.906.31"
$ws.Cells.Item(2, 5).Value = "  -1.92%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.415.30"
$ws.Cells.Item(3, 5).Value = "  -1.45%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "574.91"
$ws.Cells.Item(5, 5).Value = "  -0.79%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "148.81"
$ws.Cells.Item(6, 5).Value = "  +0.59%  "
$ws.Cells.Item(7, 5).Value = "  -0.04%  "
$ws.Cells.Item(8, 5).Value = "  +0.71%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "8.03"
$ws.Cells.Item(9, 5).Value = "  +4.62%  "
$ws.Cells.Item(10, 5).Value = "  -1.38%  "
$ws.Cells.Item(11, 5).Value = "  +1.83%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "3.990.65"
$ws.Cells.Item(12, 5).Value = "  -1.68%  "
$ws.Cells.Item(13, 5).Value = "  +0.05%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "28.23"
$ws.Cells.Item(14, 5).Value = "  -4.88%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.419.80"
$ws.Cells.Item(15, 5).Value = "  -1.13%  "
$ws.Cells.Item(16, 5).Value = "  -0.72%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "61.941.54"
$ws.Cells.Item(17, 5).Value = "  -1.81%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "6.35"
$ws.Cells.Item(18, 5).Value = "  +0.28%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "14.43"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "8.88"
$ws.Cells.Item(20, 5).Value = "  -4.31%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "380.47"
$ws.Cells.Item(21, 5).Value = "  -2.15%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.567"
$ws.Cells.Item(22, 5).Value = "  +1.15%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "75.14"
$ws.Cells.Item(24, 5).Value = "  +0.09%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.558.07"
$ws.Cells.Item(25, 5).Value = "  -1.42%  "
$ws.Cells.Item(27, 5).Value = "  +0.59%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.63"
$ws.Cells.Item(28, 5).Value = "  +0.51%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.997"
$ws.Cells.Item(29, 5).Value = "  -0.19%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "7.91"
$ws.Cells.Item(30, 5).Value = "  -3.25%  "
$ws.Cells.Item(31, 5).Value = "  -0.78%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.999"
$ws.Cells.Item(32, 5).Value = "  -0.10%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.34"
$ws.Cells.Item(33, 5).Value = "  -2.52%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "23.04"
$ws.Cells.Item(34, 5).Value = "  -2.43%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.48"
$ws.Cells.Item(35, 5).Value = "  +3.18%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.59"
$ws.Cells.Item(36, 5).Value = "  +2.25%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "169.71"
$ws.Cells.Item(37, 5).Value = "  -0.36%  "
$ws.Cells.Item(38, 5).Value = "  -2.55%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "31.01"
$ws.Cells.Item(39, 5).Value = "  -3.15%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.449.34"
$ws.Cells.Item(40, 5).Value = "  -1.54%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0773"
$ws.Cells.Item(41, 5).Value = "  +1.39%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "42.55"
$ws.Cells.Item(42, 5).Value = "  +0.56%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.775"
$ws.Cells.Item(43, 5).Value = "  -2.88%  "
$ws.Cells.Item(44, 5).Value = "  -1.46%  "
$ws.Cells.Item(45, 5).Value = "  -3.24%  "
$ws.Cells.Item(46, 5).Value = "  -4.70%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.541.96"
$ws.Cells.Item(47, 5).Value = "  -3.19%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "6.90"
$ws.Cells.Item(48, 5).Value = "  +1.91%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "22.49"
$ws.Cells.Item(49, 5).Value = "  -2.43%  "
$ws.Cells.Item(50, 5).Value = "  +0.05%  "
$ws.Cells.Item(51, 5).Value = "  -4.83%  "
